$d = $word.ActiveDocument

# 1. Extend the skills line with two more technologies.
$d.Content.Find.Execute(
    "mongodb, UML", $true, $false, $false, $false, $false,
    $true, 1, $false, "mongodb, UML, ExpressJS, MEAN stack", 2) | Out-Null

# 2. Add a leading bullet ("* ") to the Actonate internship heading.
$d.Content.Find.Execute(
    "Web Developer Intern : Actonate", $true, $false, $false, $false, $false,
    $true, 1, $false, "* Web Developer Intern : Actonate", 2) | Out-Null

# 3. Add a leading bullet ("*") to the Logic Bits heading.
$d.Content.Find.Execute(
    "Assistant Developer (Individual Contributor): Logic Bits", $true, $false, $false, $false, $false,
    $true, 1, $false, "*Assistant Developer (Individual Contributor): Logic Bits", 2) | Out-Null

# 4. Insert a leading "*" bullet run before the Veejansh Inc heading.
$r = $d.Content
$found = $r.Find.Execute(
    "Developer (Intern): Veejansh Inc", $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0)
if ($found) {
    $r.InsertBefore("*")
}
